$d = $word.ActiveDocument

# --- Paragraph 1 formatting: add a paragraph border with 5-twip spacing on
#     all sides and change the left indent from 120 to 225 twips (6pt ->
#     11.25pt), matching the formatting already used by the paragraphs that
#     follow it in the body of the document.
$p1 = $d.Paragraphs.Item(1)
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5
$p1.Format.LeftIndent = 11.25

# --- Paragraph 1 text: rename the placeholder bookmark id and drop the
#     extra run that only contained a single trailing space character.
$newId = "**ID__AFFARS_SMC_PGI_5343_102_90__ID**"

$findRange = $d.Content
$found = $findRange.Find.Execute("**ID__AFFARS_pgi_5343_topic_5__ID**", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $idStart = $findRange.Start
    $findRange.Text = $newId

    $idEnd = $idStart + $newId.Length
    $spaceRange = $d.Range($idEnd, $idEnd + 1)
    if ($spaceRange.Text -eq " ") {
        $spaceRange.Text = ""
    }
}
